$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("F4").Value = "Euclides-Usinagem"
$ws.Range("F6").Value = "Euclides-Usinagem"
